$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "ToDo" -> "AMAZING" for rows 24, 33, 38 (assign this new shared string first)
$ws.Range("N24").Value = "AMAZING"
$ws.Range("N33").Value = "AMAZING"
$ws.Range("N38").Value = "AMAZING"

# Add new "AMAZING" value to row 40 (previously empty)
$ws.Range("N40").Value = "AMAZING"

# Update the "ToDo" -> "ToDo?" for rows 10 and 14 (column N)
$ws.Range("N10").Value = "ToDo?"
$ws.Range("N14").Value = "ToDo?"

# Clear the "ToDo" value in row 32 (cell becomes empty, but keeps style)
$ws.Range("N32").Value = $null

# Update selection / view state
$ws.Range("N15").Select()
